# Weekly price update: a new price record (week) is inserted as row 7,
# pushing the existing rows 7-42 down to 8-43 (all their data is preserved
# unchanged, just shifted down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 7, shifting rows 7..42 -> 8..43
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = "Vega Modelo de Temuco"
$ws.Range("C7").Value = "La Araucanía"
$ws.Range("D7").Value = 44811
$ws.Range("E7").Value = 9
$ws.Range("F7").Value = 100112042
$ws.Range("G7").Value = "Locoto"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 2700
$ws.Range("L7").Value = 2700
$ws.Range("M7").Value = 2700
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 2700
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
